$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ProductLoanInput")
$ws2 = $wb.Worksheets.Item("ProductLoanOutput")

# Rename the product code value (drop the space before "Repayment") on both
# sheets - this is the same shared string used for B1 on both worksheets.
$newName = "968-MS-EI-DB-DL-REC-NON-RNI-CTRFD-DL-MD-TR-1-LateRepayment"
$ws1.Range("B1").Value = $newName
$ws2.Range("B1").Value = $newName

# Move the selection / active sheet so ProductLoanOutput (sheet 2) becomes
# the active tab, with B1 selected on both sheets (matches the updated
# sheetViews / workbookView activeTab in the target workbook).
$ws1.Range("B1").Select()
$ws2.Activate()
$ws2.Range("B1").Select()
